# exp_0412-1710: caps: two sentence for one images
#
# Add two new experiment rows (exp_0411-2200 and exp_0412-1710) and reword a
# couple of the existing caps notes. Cell writes below are deliberately
# ordered so that newly-introduced text lands in the shared-string table in
# the same sequence the workbook author produced it in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 (new): exp_0411-2200, caps = 4 blocks, everything else same ---
$ws.Range("A12").Value = "exp_0411-2200/checkpoint-44280"
$ws.Range("B12").Value = 71.05
$ws.Range("C12").Value = 54.19
$ws.Range("D12").Value = 40.53
$ws.Range("E12").Value = 30.29
$ws.Range("F12").Value = 25.15
$ws.Range("G12").Value = 52.55
$ws.Range("H12").Value = 98.48
$ws.Range("I12").Value = 18.4
$ws.Range("J12").Value = 53.49
$ws.Range("K12").Value = "4卡"
$ws.Range("L12").Value = "caps为4个block的句子，其他和以上一样"

# --- Row 11: reword the caps note ("caps为9个句子" -> "caps为9个block句子") ---
$ws.Range("A11").Value = "exp_0410-2134/checkpoint-44280"
$ws.Range("L11").Value = "9格 10轮 CDN，映射层一样，caps为9个block句子，bs为128，以上的bs都是256，由于句子过长，encoder的输入文本长度为293，而以上的输入文本长度为143"

# --- Row 10: reword the caps note ("caps数为9" -> "caps数为9block句子") ---
$ws.Range("L10").Value = "9格 10轮 CDN，CDN映射层修改，caps数为9block句子，以上的caps都是4"

# --- Row 13 (new): exp_0412-1710, caps = whole-image retrieval, 2 sentences ---
$ws.Range("L13").Value = "caps为整张图片检索的2个句子，其他一样"
$ws.Range("A13").Value = "exp_0412-1710/checkpoint-22140"
$ws.Range("B13").Value = 75.15
$ws.Range("C13").Value = 58.99
$ws.Range("D13").Value = 44.99
$ws.Range("E13").Value = 34.14
$ws.Range("F13").Value = 26.7
$ws.Range("G13").Value = 55.33
$ws.Range("H13").Value = 111.26
$ws.Range("I13").Value = 19.87
$ws.Range("J13").Value = 55.69
$ws.Range("K13").Value = "4卡"

# Match the numeric display format used by the other data rows (B:J, style 1
# -> numFmt "0.00_ ").
$ws.Range("B12:J13").NumberFormat = "0.00_ "

# The sheet view no longer freezes/scrolls to B1, and the last selection
# moves down two rows (L14 -> L17) to sit just below the newly appended data.
$ws.Range("L17").Select()
